# This workbook holds one weekly record (row) per date for "Acelga" prices
# at "Macroferia Regional de Talca". The commit adds a new weekly record.
#
# The new record is inserted as row 44 (pushing the previous rows 44-185
# down to become rows 45-186), and is populated with a new date (D) and a
# new volume (J), while the remaining columns reuse the same values that
# were already present on the (now shifted) row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 44; existing rows 44:185 shift to 45:186
$ws.Rows.Item(44).Insert()

# Populate the newly inserted row 44 with the new weekly record
$ws.Cells.Item(44, 1).Value  = 5
$ws.Cells.Item(44, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(44, 3).Value  = "Maule"
$ws.Cells.Item(44, 4).Value  = 44487
$ws.Cells.Item(44, 5).Value  = 7
$ws.Cells.Item(44, 6).Value  = 100112009
$ws.Cells.Item(44, 7).Value  = "Acelga"
$ws.Cells.Item(44, 8).Value  = "Sin especificar"
$ws.Cells.Item(44, 9).Value  = "Primera"
$ws.Cells.Item(44, 10).Value = 400
$ws.Cells.Item(44, 11).Value = 1800
$ws.Cells.Item(44, 12).Value = 1800
$ws.Cells.Item(44, 13).Value = 1800
$ws.Cells.Item(44, 14).Value = "$/docena de atados (4 kilos)"
$ws.Cells.Item(44, 15).Value = "Región del Maule"
$ws.Cells.Item(44, 16).Value = 450
$ws.Cells.Item(44, 17).Value = 4
$ws.Cells.Item(44, 18).Value = "Hortaliza"
